# Table_S5.xlsx revision2 edit
# - rename "C_i" sheet to "C_I"
# - rename the "NC.disease.gene.analysis" header on the "R" sheet to "disease gene prediction"
# - bold the header row + first column of the "R" sheet, add column widths / zoom
# - insert a brand-new "Legend" sheet as the first (active) tab, describing each
#   sheet of the workbook

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename "C_i" -> "C_I"
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("C_i").Name = "C_I"

# ---------------------------------------------------------------------------
# 2. Update the "R" sheet: header text, bold formatting, column widths, zoom
# ---------------------------------------------------------------------------
$rSheet = $wb.Worksheets.Item("R")
$rSheet.Range("F1").Value = "disease gene prediction"

# Bold the header row and the first (label) column
$rSheet.Range("B1:G1").Font.Bold = $true
$rSheet.Range("A2:A20").Font.Bold = $true

# Column widths (approximate "best fit" widths from the authored workbook)
$rSheet.Columns.Item(1).ColumnWidth = 6.16
$rSheet.Columns.Item(2).ColumnWidth = 26.66
$rSheet.Columns.Item(3).ColumnWidth = 8.33
$rSheet.Columns.Item(4).ColumnWidth = 9.33
$rSheet.Columns.Item(5).ColumnWidth = 14.83
$rSheet.Columns.Item(6).ColumnWidth = 20
$rSheet.Columns.Item(7).ColumnWidth = 5.16

$rSheet.Activate()
$rSheet.Range("F2").Select()
$excel.ActiveWindow.Zoom = 142

# ---------------------------------------------------------------------------
# 3. Insert the new "Legend" sheet before the "R" sheet (i.e. as the first tab)
# ---------------------------------------------------------------------------
$firstSheet = $wb.Worksheets.Item(1)
$legend = $wb.Worksheets.Add($firstSheet)
$legend.Name = "Legend"

$legendData = @(
    @("R", "aggregate correlation"),
    @("C_I", "weight matrix"),
    @("C_T", "correlation matrix relative to topological analyses"),
    @("C_PC", "correlation matrix relative to pathway complexes"),
    @("C_PATH", "correlation matrix relative to pathways"),
    @("C_PCT", "correlation matrix relative to pathway cross-talk"),
    @("C_DGP", "correlation matrix relative to disease gene prioritization")
)

$r = 1
foreach ($row in $legendData) {
    $legend.Cells.Item($r, 1).Value = $row[0]
    $legend.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}

$legend.Columns.Item(2).ColumnWidth = 39.33

$legend.Activate()
$legend.Range("B8").Select()
$excel.ActiveWindow.Zoom = 186
